$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 777/778, shifting the existing rows 777:834 down to 779:836.
$ws.Rows("777:778").Insert()

# Row 777 - new weekly entry (Escarola, Primera)
$ws.Range("A777").Value = 4
$ws.Range("B777").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C777").Value = "Los Lagos"
$ws.Range("D777").Value = 45013
$ws.Range("E777").Value = 10
$ws.Range("F777").Value = 100112033
$ws.Range("G777").Value = "Lechuga"
$ws.Range("H777").Value = "Escarola"
$ws.Range("I777").Value = "Primera"
$ws.Range("J777").Value = 300
$ws.Range("K777").Value = 11000
$ws.Range("L777").Value = 11000
$ws.Range("M777").Value = 11000
$ws.Range("N777").Value = "$/caja 15 unidades"
$ws.Range("O777").Value = "Región de Coquimbo"
$ws.Range("P777").Value = 733
$ws.Range("Q777").Value = 15
$ws.Range("R777").Value = "Hortaliza"

# Row 778 - new weekly entry (Escarola, Segunda)
$ws.Range("A778").Value = 4
$ws.Range("B778").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C778").Value = "Los Lagos"
$ws.Range("D778").Value = 45013
$ws.Range("E778").Value = 10
$ws.Range("F778").Value = 100112033
$ws.Range("G778").Value = "Lechuga"
$ws.Range("H778").Value = "Escarola"
$ws.Range("I778").Value = "Segunda"
$ws.Range("J778").Value = 300
$ws.Range("K778").Value = 10000
$ws.Range("L778").Value = 10000
$ws.Range("M778").Value = 10000
$ws.Range("N778").Value = "$/caja 18 unidades"
$ws.Range("O778").Value = "Región de Coquimbo"
$ws.Range("P778").Value = 556
$ws.Range("Q778").Value = 18
$ws.Range("R778").Value = "Hortaliza"
